{"js": "// Update the title date line and the 25 division-problem answers in the\n// table. Cell values are addressed positionally (row, col) because several\n// answer strings repeat (e.g. \"66\u00f77=9, 3\" and \"82\u00f79=9, 1\" each occur twice),\n// so a text search-and-replace-all would corrupt the document.\n\n// 1) Title paragraph: \"2026-02-11 Wednesday\" -> \"2026-02-12 Thursday\"\nconst body = context.document.body;\nconst titleHits = body.search(\"2026-02-11 Wednesday\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"2026-02-12 Thursday\", \"Replace\");\n}\n\n// 2) Table cell answers, addressed by (row, col) within the single table.\nconst table = context.document.body.tables.getFirst();\n\nconst updates = [\n  { row: 0, col: 0, value: \"37\u00f77=5, 2\" },\n  { row: 0, col: 1, value: \"38\u00f76=6, 2\" },\n  { row: 0, col: 2, value: \"69\u00f72=34, 1\" },\n  { row: 0, col: 3, value: \"27\u00f73=9, 0\" },\n  { row: 0, col: 4, value: \"97\u00f73=32, 1\" },\n\n  { row: 4, col: 0, value: \"23\u00f75=4, 3\" },\n  { row: 4, col: 1, value: \"94\u00f79=10, 4\" },\n  { row: 4, col: 2, value: \"87\u00f75=17, 2\" },\n  { row: 4, col: 3, value: \"90\u00f74=22, 2\" },\n  { row: 4, col: 4, value: \"24\u00f73=8, 0\" },\n\n  { row: 8, col: 0, value: \"72\u00f77=10, 2\" },\n  { row: 8, col: 1, value: \"59\u00f78=7, 3\" },\n  { row: 8, col: 2, value: \"55\u00f74=13, 3\" },\n  { row: 8, col: 3, value: \"83\u00f77=11, 6\" },\n  { row: 8, col: 4, value: \"93\u00f77=13, 2\" },\n\n  { row: 12, col: 0, value: \"74\u00f78=9, 2\" },\n  { row: 12, col: 1, value: \"22\u00f74=5, 2\" },\n  { row: 12, col: 2, value: \"12\u00f76=2, 0\" },\n  { row: 12, col: 3, value: \"34\u00f75=6, 4\" },\n  { row: 12, col: 4, value: \"90\u00f77=12, 6\" },\n\n  { row: 16, col: 0, value: \"99\u00f77=14, 1\" },\n  { row: 16, col: 1, value: \"10\u00f72=5, 0\" },\n  { row: 16, col: 2, value: \"37\u00f75=7, 2\" },\n  { row: 16, col: 3, value: \"60\u00f76=10, 0\" },\n  { row: 16, col: 4, value: \"43\u00f75=8, 3\" }\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.value = u.value;\n}\n\nawait context.sync();\n", "ps1": "# Update the title date line and the 25 division-problem answers in the\n# table. Cell values are addressed positionally (row, col) because several\n# answer strings repeat (e.g. \"66\u00f77=9, 3\" and \"82\u00f79=9, 1\" each occur twice),\n# so a blind FindReplace-all would corrupt the document.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2026-02-11 Wednesday\" -> \"2026-02-12 Thursday\"\n$find = $d.Content.Find\n$find.Text = \"2026-02-11 Wednesday\"\n$find.Replacement.Text = \"2026-02-12 Thursday\"\n$find.Execute(\n    \"2026-02-11 Wednesday\", $false, $false, $false, $false, $false,\n    $true, 1, $false, \"2026-02-12 Thursday\", 2\n) | Out-Null\n\n# 2) Table cell answers, addressed by 1-based (row, col) within the table.\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text  = \"37\u00f77=5, 2\"\n$t.Cell(1, 2).Range.Text  = \"38\u00f76=6, 2\"\n$t.Cell(1, 3).Range.Text  = \"69\u00f72=34, 1\"\n$t.Cell(1, 4).Range.Text  = \"27\u00f73=9, 0\"\n$t.Cell(1, 5).Range.Text  = \"97\u00f73=32, 1\"\n\n$t.Cell(5, 1).Range.Text  = \"23\u00f75=4, 3\"\n$t.Cell(5, 2).Range.Text  = \"94\u00f79=10, 4\"\n$t.Cell(5, 3).Range.Text  = \"87\u00f75=17, 2\"\n$t.Cell(5, 4).Range.Text  = \"90\u00f74=22, 2\"\n$t.Cell(5, 5).Range.Text  = \"24\u00f73=8, 0\"\n\n$t.Cell(9, 1).Range.Text  = \"72\u00f77=10, 2\"\n$t.Cell(9, 2).Range.Text  = \"59\u00f78=7, 3\"\n$t.Cell(9, 3).Range.Text  = \"55\u00f74=13, 3\"\n$t.Cell(9, 4).Range.Text  = \"83\u00f77=11, 6\"\n$t.Cell(9, 5).Range.Text  = \"93\u00f77=13, 2\"\n\n$t.Cell(13, 1).Range.Text = \"74\u00f78=9, 2\"\n$t.Cell(13, 2).Range.Text = \"22\u00f74=5, 2\"\n$t.Cell(13, 3).Range.Text = \"12\u00f76=2, 0\"\n$t.Cell(13, 4).Range.Text = \"34\u00f75=6, 4\"\n$t.Cell(13, 5).Range.Text = \"90\u00f77=12, 6\"\n\n$t.Cell(17, 1).Range.Text = \"99\u00f77=14, 1\"\n$t.Cell(17, 2).Range.Text = \"10\u00f72=5, 0\"\n$t.Cell(17, 3).Range.Text = \"37\u00f75=7, 2\"\n$t.Cell(17, 4).Range.Text = \"60\u00f76=10, 0\"\n$t.Cell(17, 5).Range.Text = \"43\u00f75=8, 3\"\n"}
